# Updates the "cryptos" list with refreshed price/volume figures
# (and a few re-ranked rows) as produced by the scheduled GitHub
# Actions scrape on Fri Sep  8 07:23:52 UTC 2023.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Addr,
        [string]$Text
    )
    $rng = $ws.Range($Addr)
    # Force the cell to remain plain text even when the string looks
    # like a number (e.g. "4.30", "0.418"), then drop back to the
    # workbook's default "Normal" style so no stray formatting is left
    # behind on the cell.
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "26.262.52"
Set-TextValue "E2" "  +1.84%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.646.96"
Set-TextValue "E3" "  +0.51%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  -0.16%  "

# Row 5 - BNB
Set-TextValue "D5" "216.97"
Set-TextValue "E5" "  +0.60%  "

# Row 6 - XRP
Set-TextValue "D6" "0.505"
Set-TextValue "E6" "  +0.49%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.259"
Set-TextValue "E8" "  +0.27%  "

# Row 9 - Dogecoin
Set-TextValue "E9" "  +0.42%  "

# Row 10 - Solana
Set-TextValue "D10" "19.95"
Set-TextValue "E10" "  +1.49%  "

# Row 11 - TRON
Set-TextValue "D11" "0.0793"
Set-TextValue "E11" "  +0.04%  "

# Row 12 - was WrappedliquidstakedEther2.0, now WrappedEther
Set-TextValue "B12" "WrappedEther"
Set-TextValue "C12" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D12" "1.718.35"
Set-TextValue "E12" "  +4.84%  "

# Row 13 - was Polkadot, now WrappedliquidstakedEther2.0
Set-TextValue "B13" "WrappedliquidstakedEther2.0"
Set-TextValue "C13" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D13" "1.873.70"
Set-TextValue "E13" "  +0.46%  "

# Row 14 - was WrappedEther, now Polkadot
Set-TextValue "B14" "Polkadot"
Set-TextValue "C14" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D14" "4.30"
Set-TextValue "E14" "  +0.92%  "

# Row 15 - Polygon
Set-TextValue "E15" "  -2.74%  "

# Row 16 - ShibaInu
Set-TextValue "E16" "  +0.30%  "

# Row 17 - Litecoin
Set-TextValue "D17" "63.35"
Set-TextValue "E17" "  +0.32%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "26.259.87"
Set-TextValue "E18" "  +1.66%  "

# Row 19 - Dai
Set-TextValue "E19" "  -0.13%  "

# Row 20 - Uniswap
Set-TextValue "D20" "4.45"
Set-TextValue "E20" "  -0.62%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "195.78"
Set-TextValue "E21" "  +1.74%  "

# Row 23 - Chainlink
Set-TextValue "E23" "  +0.33%  "

# Row 24 - Toncoin
Set-TextValue "E24" "  -3.59%  "

# Row 25 - Monero
Set-TextValue "D25" "143.42"
Set-TextValue "E25" "  +1.07%  "

# Row 26 - BinanceUSD
Set-TextValue "E26" "  -0.20%  "

# Row 27 - Stellar
Set-TextValue "E27" "  +0.90%  "

# Row 28 - Cosmos
Set-TextValue "E28" "  +0.42%  "

# Row 29 - EthereumClassic
Set-TextValue "D29" "15.66"
Set-TextValue "E29" "  +1.00%  "

# Row 30 - PancakeSwap
Set-TextValue "E30" "  +1.01%  "

# Row 31 - Hedera
Set-TextValue "D31" "0.0504"
Set-TextValue "E31" "  +2.35%  "

# Row 32 - InternetComputer(DFINITY)
Set-TextValue "E32" "  +0.96%  "

# Row 33 - Filecoin
Set-TextValue "D33" "3.26"
Set-TextValue "E33" "  +0.95%  "

# Row 34 - LidoDAOToken
Set-TextValue "E34" "  +2.04%  "

# Row 35 - HuobiToken
Set-TextValue "E35" "  +1.24%  "

# Row 37 - Maker
Set-TextValue "D37" "1.139.14"
Set-TextValue "E37" "  +0.34%  "

# Row 38 - ImmutableX
Set-TextValue "D38" "0.554"
Set-TextValue "E38" "  +1.55%  "

# Row 39 - MXToken
Set-TextValue "E39" "  -1.73%  "

# Row 40 - VeChain
Set-TextValue "E40" "  +1.23%  "

# Row 41 - PaxDollar
Set-TextValue "E41" "  -0.17%  "

# Row 42 - Quant
Set-TextValue "D42" "100.55"
Set-TextValue "E42" "  +0.00%  "

# Row 43 - FraxShare
Set-TextValue "D43" "5.52"
Set-TextValue "E43" "  -1.00%  "

# Row 44 - TrustWalletToken
Set-TextValue "D44" "0.803"
Set-TextValue "E44" "  -0.15%  "

# Row 45 - RocketPoolETH
Set-TextValue "D45" "1.782.56"
Set-TextValue "E45" "  +0.46%  "

# Row 46 - Aave
Set-TextValue "D46" "56.98"
Set-TextValue "E46" "  +3.05%  "

# Row 47 - RenderToken
Set-TextValue "E47" "  +4.14%  "

# Row 48 - Cronos
Set-TextValue "E48" "  +3.02%  "

# Row 49 - was EnergySwap, now Mantle
Set-TextValue "B49" "Mantle"
Set-TextValue "C49" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D49" "0.418"
Set-TextValue "E49" "  +0.19%  "

# Row 50 - was Mantle, now EnergySwap
Set-TextValue "B50" "EnergySwap"
Set-TextValue "C50" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D50" "7.69"
Set-TextValue "E50" "  +3.20%  "

# Row 51 - Algorand
Set-TextValue "E51" "  +1.95%  "
